$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.802135104405522
$ws.Range("C2").Value = 0.3891317079198302
$ws.Range("D2").Value = 0.6669206992951331
$ws.Range("E2").Value = 0.2718350268723029
$ws.Range("G2").Value = 2.434755334442031
$ws.Range("H2").Value = 1.886318112648723
$ws.Range("J2").Value = 0.1413281424595567

$ws.Range("B3").Value = 1.689558433147681
$ws.Range("C3").Value = 0.3582393212594468
$ws.Range("D3").Value = 0.6576310717023546
$ws.Range("E3").Value = 0.2671919136001009
$ws.Range("G3").Value = 2.373819228495535
$ws.Range("H3").Value = 1.865367232599567
$ws.Range("J3").Value = 0.1381904177401978

$ws.Range("B4").Value = 1.621589877640758
$ws.Range("C4").Value = 0.3395084973742541
$ws.Range("D4").Value = 0.6523081617365278
$ws.Range("E4").Value = 0.2645036094725768
$ws.Range("G4").Value = 2.337962502499778
$ws.Range("H4").Value = 1.853537980294675
$ws.Range("J4").Value = 0.1363535972871119

$ws.Range("B5").Value = 1.594180636723536
$ws.Range("C5").Value = 0.3319346987271956
$ws.Range("D5").Value = 0.6502344270770948
$ws.Range("E5").Value = 0.2634487776029744
$ws.Range("G5").Value = 2.323737897705001
$ws.Range("H5").Value = 1.848975667595312
$ws.Range("J5").Value = 0.1356274987801669

$ws.Range("B6").Value = 1.58964673008569
$ws.Range("C6").Value = 0.3306806347700899
$ws.Range("D6").Value = 0.6498958350793487
$ws.Range("E6").Value = 0.2632760741359377
$ws.Range("G6").Value = 2.321399172811141
$ws.Range("H6").Value = 1.848233636512987
$ws.Range("J6").Value = 0.135508280534971

$ws.Range("B7").Value = 1.621219061077625
$ws.Range("C7").Value = 0.339406115551725
$ws.Range("D7").Value = 0.6522798088965942
$ws.Range("E7").Value = 0.2644892192270163
$ws.Range("G7").Value = 2.337769102304833
$ws.Range("H7").Value = 1.853475408482353
$ws.Range("J7").Value = 0.1363437142823969

$ws.Range("B8").Value = 1.763078089718306
$ws.Range("C8").Value = 0.3784304662051454
$ws.Range("D8").Value = 0.6636383191798814
$ws.Range("E8").Value = 0.270200207070495
$ws.Range("G8").Value = 2.413418377106524
$ws.Range("H8").Value = 1.878878300684761
$ws.Range("J8").Value = 0.1402275345665771

$ws.Range("B9").Value = 2.050514955470817
$ws.Range("C9").Value = 0.4568666952109766
$ws.Range("D9").Value = 0.6889554992229137
$ws.Range("E9").Value = 0.282700029089888
$ws.Range("G9").Value = 2.574345834378363
$ws.Range("H9").Value = 1.936996759959186
$ws.Range("J9").Value = 0.148563244009928

$ws.Range("B10").Value = 2.267495682378581
$ws.Range("C10").Value = 0.5157034099338489
$ws.Range("D10").Value = 0.7094429611132966
$ws.Range("E10").Value = 0.2926927000004582
$ws.Range("G10").Value = 2.700560783556853
$ws.Range("H10").Value = 1.984895272953707
$ws.Range("J10").Value = 0.1551373350398109

$ws.Range("B11").Value = 2.367503136521293
$ws.Range("C11").Value = 0.5427430697430964
$ws.Range("D11").Value = 0.7191802061817896
$ws.Range("E11").Value = 0.2974179737794529
$ws.Range("G11").Value = 2.759781383129791
$ws.Range("H11").Value = 2.00784473781971
$ws.Range("J11").Value = 0.1582282642519601

$ws.Range("B12").Value = 2.405563231881843
$ws.Range("C12").Value = 0.5530225316312567
$ws.Range("D12").Value = 0.7229280182389175
$ws.Range("E12").Value = 0.2992334126267053
$ws.Range("G12").Value = 2.782471618890554
$ws.Range("H12").Value = 2.016704317588051
$ws.Range("J12").Value = 0.159413341989989

$ws.Range("B13").Value = 2.397357847148214
$ws.Range("C13").Value = 0.5508068701616935
$ws.Range("D13").Value = 0.7221181595958512
$ws.Range("E13").Value = 0.2988412612537132
$ws.Range("G13").Value = 2.777573016154577
$ws.Range("H13").Value = 2.014788690425291
$ws.Range("J13").Value = 0.1591574615914482

$ws.Range("B14").Value = 2.370630553584476
$ws.Range("C14").Value = 0.5435879591767048
$ws.Range("D14").Value = 0.719487324892981
$ws.Range("E14").Value = 0.2975668068263388
$ws.Range("G14").Value = 2.761642788146276
$ws.Range("H14").Value = 2.008570218272666
$ws.Range("J14").Value = 0.1583254674547021

$ws.Range("B15").Value = 2.354284053713457
$ws.Range("C15").Value = 0.539171413596705
$ws.Range("D15").Value = 0.717883760599392
$ws.Range("E15").Value = 0.2967895709356867
$ws.Range("G15").Value = 2.751919686738972
$ws.Range("H15").Value = 2.004783316606961
$ws.Range("J15").Value = 0.1578177551057394

$ws.Range("B16").Value = 2.260986345689162
$ws.Range("C16").Value = 0.5139418840631151
$ws.Range("D16").Value = 0.7088150514007339
$ws.Range("E16").Value = 0.2923875274860208
$ws.Range("G16").Value = 2.69672733410161
$ws.Range("H16").Value = 1.983419011755899
$ws.Range("J16").Value = 0.1549373700526928

$ws.Range("B17").Value = 2.20408621931972
$ws.Range("C17").Value = 0.4985351560888489
$ws.Range("D17").Value = 0.7033589539820184
$ws.Range("E17").Value = 0.2897331884028773
$ws.Range("G17").Value = 2.663334509130152
$ws.Range("H17").Value = 1.970611468691146
$ws.Range("J17").Value = 0.1531961745293131

$ws.Range("B18").Value = 2.171481166538229
$ws.Range("C18").Value = 0.4896994143367124
$ws.Range("D18").Value = 0.7002600071398604
$ws.Range("E18").Value = 0.2882233676391053
$ws.Range("G18").Value = 2.644297402415532
$ws.Range("H18").Value = 1.963353973050431
$ws.Range("J18").Value = 0.1522041185381369

$ws.Range("B19").Value = 2.160462610587956
$ws.Range("C19").Value = 0.4867121959240421
$ws.Range("D19").Value = 0.6992174840129053
$ws.Range("E19").Value = 0.2877150604469492
$ws.Range("G19").Value = 2.63788073179677
$ws.Range("H19").Value = 1.960915372619638
$ws.Range("J19").Value = 0.1518698402167473

$ws.Range("B20").Value = 2.210130651561769
$ws.Range("C20").Value = 0.5001725536866388
$ws.Range("D20").Value = 0.703935698378757
$ws.Range("E20").Value = 0.2900139979468292
$ws.Range("G20").Value = 2.666871646105079
$ws.Range("H20").Value = 1.971963549344849
$ws.Range("J20").Value = 0.1533805501246661

$ws.Range("B21").Value = 2.378475847830543
$ws.Range("C21").Value = 0.5457072338941771
$ws.Range("D21").Value = 0.720258417971479
$ws.Range("E21").Value = 0.2979404350702737
$ws.Range("G21").Value = 2.766314657527857
$ws.Range("H21").Value = 2.010392126191789
$ws.Range("J21").Value = 0.1585694463296932

$ws.Range("B22").Value = 2.489604490846261
$ws.Range("C22").Value = 0.5757010073887727
$ws.Range("D22").Value = 0.7312793118694287
$ws.Range("E22").Value = 0.3032729608342635
$ws.Range("G22").Value = 2.832851610610646
$ws.Range("H22").Value = 2.036494239094282
$ws.Range("J22").Value = 0.162045934736355

$ws.Range("B23").Value = 2.430191140319948
$ws.Range("C23").Value = 0.5596711118884627
$ws.Range("D23").Value = 0.725364777621337
$ws.Range("E23").Value = 0.3004128829294856
$ws.Range("G23").Value = 2.797196480808338
$ws.Range("H23").Value = 2.022471994147566
$ws.Range("J23").Value = 0.1601826066966368

$ws.Range("B24").Value = 2.207397628801061
$ws.Range("C24").Value = 0.4994322184062412
$ws.Range("D24").Value = 0.7036748344233388
$ws.Range("E24").Value = 0.2898869935305086
$ws.Range("G24").Value = 2.665272005783748
$ws.Range("H24").Value = 1.97135194452332
$ws.Range("J24").Value = 0.1532971659593869

$ws.Range("B25").Value = 1.971748911304303
$ws.Range("C25").Value = 0.4354385768227758
$ws.Range("D25").Value = 0.6817775049526631
$ws.Range("E25").Value = 0.279177608410599
$ws.Range("G25").Value = 2.529431850769811
$ws.Range("H25").Value = 1.920371187916118
$ws.Range("J25").Value = 0.1413281424595567
